$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6160493638826665
$ws.Range("C2").Value = 0.138989126959995
$ws.Range("E2").Value = 0.2167433957142926
$ws.Range("F2").Value = 2.100022427693361
$ws.Range("G2").Value = 0.002482084203086887
$ws.Range("I2").Value = 1.031830816547593
$ws.Range("J2").Value = 0.06394189639972581
$ws.Range("K2").Value = 0.3300181322398146
$ws.Range("L2").Value = 0.4664012796701797
$ws.Range("O2").Value = 3.873840421751538
$ws.Range("B3").Value = 0.5725389478365059
$ws.Range("C3").Value = 0.1386822754600878
$ws.Range("E3").Value = 0.215328455227084
$ws.Range("F3").Value = 2.102601023592925
$ws.Range("G3").Value = 0.002484395431093864
$ws.Range("I3").Value = 1.041599121772443
$ws.Range("J3").Value = 0.06270644475393539
$ws.Range("K3").Value = 0.2938163420836304
$ws.Range("L3").Value = 0.4557100808996921
$ws.Range("O3").Value = 3.907300601650846
$ws.Range("B4").Value = 0.5459411450598566
$ws.Range("C4").Value = 0.1385004160860568
$ws.Range("E4").Value = 0.2145390260221056
$ws.Range("F4").Value = 2.105214164293919
$ws.Range("G4").Value = 0.002485890243645423
$ws.Range("I4").Value = 1.048074538806837
$ws.Range("J4").Value = 0.06194462078331853
$ws.Range("K4").Value = 0.2715765418725624
$ws.Range("L4").Value = 0.4493236593938832
$ws.Range("O4").Value = 3.929784546912003
$ws.Range("B5").Value = 0.5351327306454152
$ws.Range("C5").Value = 0.1384279696707154
$ws.Range("E5").Value = 0.2142373480581021
$ws.Range("F5").Value = 2.106538267356349
$ws.Range("G5").Value = 0.002486518485060215
$ws.Range("I5").Value = 1.050833458075804
$ws.Range("J5").Value = 0.06163337388616341
$ws.Range("K5").Value = 0.2625112601028832
$ws.Range("L5").Value = 0.4467661411850514
$ws.Range("O5").Value = 3.93943462706109
$ws.Range("B6").Value = 0.533339861495449
$ws.Range("C6").Value = 0.1384160409004274
$ws.Range("E6").Value = 0.214188465861703
$ws.Range("F6").Value = 2.106773797898754
$ws.Range("G6").Value = 0.002486623958925354
$ws.Range("I6").Value = 1.051298831166708
$ws.Range("J6").Value = 0.06158164400686061
$ws.Range("K6").Value = 0.2610058499865886
$ws.Range("L6").Value = 0.4463441921090947
$ws.Range("O6").Value = 3.941066472201342
$ws.Range("B7").Value = 0.5457952548374863
$ws.Range("C7").Value = 0.1384994322947435
$ws.Range("E7").Value = 0.2145348763340529
$ws.Range("F7").Value = 2.105230971670053
$ws.Range("G7").Value = 0.002485898638936587
$ws.Range("I7").Value = 1.048111260101141
$ws.Range("J7").Value = 0.06194042640078834
$ws.Range("K7").Value = 0.2714542931980333
$ws.Range("L7").Value = 0.4492889852893853
$ws.Range("O7").Value = 3.929912716570684
$ws.Range("B8").Value = 0.6010229900237221
$ws.Range("C8").Value = 0.1388819748167833
$ws.Range("E8").Value = 0.2162390948742718
$ws.Range("F8").Value = 2.100697927479828
$ws.Range("G8").Value = 0.002482865434685497
$ws.Range("I8").Value = 1.035099811734902
$ws.Range("J8").Value = 0.06351659931754128
$ws.Range("K8").Value = 0.31753852676826
$ws.Range("L8").Value = 0.4626781443300985
$ws.Range("O8").Value = 3.884975026731524
$ws.Range("B9").Value = 0.7102296808703272
$ws.Range("C9").Value = 0.13968345485074
$ws.Range("E9").Value = 0.2202080482539799
$ws.Range("F9").Value = 2.099971011418575
$ws.Range("G9").Value = 0.002477515510248325
$ws.Range("I9").Value = 1.013372965933879
$ws.Range("J9").Value = 0.06658092163084461
$ws.Range("K9").Value = 0.4077955681454171
$ws.Range("L9").Value = 0.4903380556640258
$ws.Range("O9").Value = 3.812238044646165
$ws.Range("B10").Value = 0.790984810748995
$ws.Range("C10").Value = 0.1403028041037544
$ws.Range("E10").Value = 0.2235032619565658
$ws.Range("F10").Value = 2.104403673382052
$ws.Range("G10").Value = 0.002473946024541602
$ws.Range("I10").Value = 0.9997180911870522
$ws.Range("J10").Value = 0.06881529724272184
$ws.Range("K10").Value = 0.4740166095862719
$ws.Range("L10").Value = 0.5115064422680717
$ws.Range("O10").Value = 3.768176677358554
$ws.Range("B11").Value = 0.8278295191639415
$ws.Range("C11").Value = 0.1405910276804434
$ws.Range("E11").Value = 0.2250840802408192
$ws.Range("F11").Value = 2.10749692096401
$ws.Range("G11").Value = 0.002472399825259214
$ws.Range("I11").Value = 0.9940069014785529
$ws.Range("J11").Value = 0.06982792800462789
$ws.Range("K11").Value = 0.5041184649205093
$ws.Range("L11").Value = 0.5213185448060926
$ws.Range("O11").Value = 3.750168383433476
$ws.Range("B12").Value = 0.8417966097582621
$ws.Range("C12").Value = 0.140701087560771
$ws.Range("E12").Value = 0.2256943978121555
$ws.Range("F12").Value = 2.108822892053425
$ws.Range("G12").Value = 0.002471825418424964
$ws.Range("I12").Value = 0.9919161716539229
$ws.Range("J12").Value = 0.07021082077405794
$ws.Range("K12").Value = 0.5155135406235161
$ws.Range("L12").Value = 0.5250601705916864
$ws.Range("O12").Value = 3.743641811547263
$ws.Range("B13").Value = 0.8387879021307469
$ws.Range("C13").Value = 0.1406773437111397
$ws.Range("E13").Value = 0.2255624359043509
$ws.Range("F13").Value = 2.10853044693917
$ws.Range("G13").Value = 0.002471948634110604
$ws.Range("I13").Value = 0.9923632465933032
$ws.Range("J13").Value = 0.0701283836068356
$ws.Range("K13").Value = 0.5130595879542454
$ws.Range("L13").Value = 0.5242531923149301
$ws.Range("O13").Value = 3.745034402854941
$ws.Range("B14").Value = 0.8289783083480131
$ws.Range("C14").Value = 0.1406000641071898
$ws.Range("E14").Value = 0.225134057419325
$ws.Range("F14").Value = 2.107602911735583
$ws.Range("G14").Value = 0.002472352346141194
$ws.Range("I14").Value = 0.9938334533860278
$ws.Range("J14").Value = 0.06985944034527591
$ws.Range("K14").Value = 0.5050560252221601
$ws.Range("L14").Value = 0.521625851039758
$ws.Range("O14").Value = 3.7496255702452
$ws.Range("B15").Value = 0.8229715511540121
$ws.Range("C15").Value = 0.1405528469268447
$ws.Range("E15").Value = 0.2248731843722709
$ws.Range("F15").Value = 2.107054899929125
$ws.Range("G15").Value = 0.002472601076293958
$ws.Range("I15").Value = 0.9947433707626949
$ws.Range("J15").Value = 0.06969463012945454
$ws.Range("K15").Value = 0.5001530930851743
$ws.Range("L15").Value = 0.5200199066854765
$ws.Range("O15").Value = 3.752475923690554
$ws.Range("B16").Value = 0.7885790393012542
$ws.Range("C16").Value = 0.1402840970292303
$ws.Range("E16").Value = 0.2234015921546089
$ws.Range("F16").Value = 2.104223171870174
$ws.Range("G16").Value = 0.00247404862901731
$ws.Range("I16").Value = 1.00010140226793
$ws.Range("J16").Value = 0.06874904110355828
$ws.Range("K16").Value = 0.4720488783089536
$ws.Range("L16").Value = 0.5108688494569691
$ws.Range("O16").Value = 3.769394527217543
$ws.Range("B17").Value = 0.7675076097976614
$ws.Range("C17").Value = 0.140120875716029
$ws.Range("E17").Value = 0.2225197228388218
$ws.Range("F17").Value = 2.10276163976657
$ws.Range("G17").Value = 0.002474956489150869
$ws.Range("I17").Value = 1.003516569953522
$ws.Range("J17").Value = 0.0681679641963342
$ws.Range("K17").Value = 0.4548016675778399
$ws.Range("L17").Value = 0.5053015381734269
$ws.Range("O17").Value = 3.780294928088011
$ws.Range("B18").Value = 0.7553981715030318
$ws.Range("C18").Value = 0.1400276058488288
$ws.Range("E18").Value = 0.2220201992564093
$ws.Range("F18").Value = 2.102022372492627
$ws.Range("G18").Value = 0.00247548597141496
$ws.Range("I18").Value = 1.00552799142671
$ws.Range("J18").Value = 0.06783338798888749
$ws.Range("K18").Value = 0.4448794538692482
$ws.Range("L18").Value = 0.5021165602746862
$ws.Range("O18").Value = 3.786756130874949
$ws.Range("B19").Value = 0.7512999169979935
$ws.Range("C19").Value = 0.1399961316476563
$ws.Range("E19").Value = 0.2218523944525153
$ws.Range("F19").Value = 2.1017894858646
$ws.Range("G19").Value = 0.002475666501217891
$ws.Range("I19").Value = 1.006217114755543
$ws.Range("J19").Value = 0.06772004574833801
$ws.Range("K19").Value = 0.4415196254196303
$ws.Range("L19").Value = 0.501041141814909
$ws.Range("O19").Value = 3.788976684703641
$ws.Range("B20").Value = 0.7697496405305344
$ws.Range("C20").Value = 0.1401381878107131
$ws.Range("E20").Value = 0.2226128024575154
$ws.Range("F20").Value = 2.102906732791368
$ws.Range("G20").Value = 0.002474859090126986
$ws.Range("I20").Value = 1.003148144051423
$ws.Range("J20").Value = 0.06822985784649305
$ws.Range("K20").Value = 0.4566378817695806
$ws.Range("L20").Value = 0.5058924104175588
$ws.Range("O20").Value = 3.779114733329436
$ws.Range("B21").Value = 0.8318592289094227
$ws.Range("C21").Value = 0.1406227382684406
$ws.Range("E21").Value = 0.2252595656809078
$ws.Range("F21").Value = 2.107871156566446
$ws.Range("G21").Value = 0.002472233465165671
$ws.Range("I21").Value = 0.9933996646234071
$ws.Range("J21").Value = 0.06993845110878283
$ws.Range("K21").Value = 0.5074069748373802
$ws.Range("L21").Value = 0.5223968613209848
$ws.Range("O21").Value = 3.748269086506866
$ws.Range("B22").Value = 0.8725373057669117
$ws.Range("C22").Value = 0.1409447512040458
$ws.Range("E22").Value = 0.2270575200791498
$ws.Range("F22").Value = 2.112016867251299
$ws.Range("G22").Value = 0.002470582173721869
$ws.Range("I22").Value = 0.9874479904741307
$ws.Range("J22").Value = 0.07105179140246065
$ws.Range("K22").Value = 0.5405648231540567
$ws.Range("L22").Value = 0.5333348809500507
$ws.Range("O22").Value = 3.72981624737605
$ws.Range("B23").Value = 0.8508190720511948
$ws.Range("C23").Value = 0.1407724041830889
$ws.Range("E23").Value = 0.2260917055004725
$ws.Range("F23").Value = 2.109721838506559
$ws.Range("G23").Value = 0.002471457594990314
$ws.Range("I23").Value = 0.9905861230664108
$ws.Range("J23").Value = 0.07045789229262311
$ws.Range("K23").Value = 0.5228701352320968
$ws.Range("L23").Value = 0.5274832820262247
$ws.Range("O23").Value = 3.739508695564439
$ws.Range("B24").Value = 0.7687360033991126
$ws.Range("C24").Value = 0.1401303592429244
$ws.Range("E24").Value = 0.2225706978738522
$ws.Range("F24").Value = 2.102840821651782
$ws.Range("G24").Value = 0.002474903100720837
$ws.Range("I24").Value = 1.003314559896392
$ws.Range("J24").Value = 0.06820187730570026
$ws.Range("K24").Value = 0.4558077497289617
$ws.Range("L24").Value = 0.5056252280193405
$ws.Range("O24").Value = 3.779647693784511
$ws.Range("B25").Value = 0.6805926682066854
$ws.Range("C25").Value = 0.1394612210977542
$ws.Range("E25").Value = 0.2190675286061712
$ws.Range("F25").Value = 2.09929477852792
$ws.Range("G25").Value = 0.002478899137153836
$ws.Range("I25").Value = 1.01884524556425
$ws.Range("J25").Value = 0.0657548658906677
$ws.Range("K25").Value = 0.3833930641216341
$ws.Range("L25").Value = 0.4827059753845901
$ws.Range("O25").Value = 3.830268289936541
